# Updated cryptos list on Mon Nov 20 21:36:18 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'37.390.39"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = "'2.024.22"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.94%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = "'253.72"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.26%  '
$ws.Range('D6').Value = "'0.620"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.57%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = "'56.92"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -6.77%  '
$ws.Range('D9').Value = "'0.385"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('D10').Value = "'57.15"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('D11').Value = "'0.0787"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.82%  '
$ws.Range('E12').Value = '  -1.83%  '
$ws.Range('D13').Value = "'14.72"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.92%  '
$ws.Range('D14').Value = "'2.323.15"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.05%  '
$ws.Range('D15').Value = "'0.815"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.74%  '
$ws.Range('D16').Value = "'21.14"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -4.72%  '
$ws.Range('D17').Value = "'5.33"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.80%  '
$ws.Range('D18').Value = "'2.035.32"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.40%  '
$ws.Range('D19').Value = "'37.244.34"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('D21').Value = "'0.0₃0847"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.87%  '
$ws.Range('D22').Value = "'5.16"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').Value = "'228.20"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.88%  '
$ws.Range('D24').Value = "'0.999"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('E25').Value = '  +3.54%  '
$ws.Range('E26').Value = '  -1.01%  '
$ws.Range('D27').Value = "'162.55"
$ws.Range('D27').ClearFormats()
$ws.Range('D28').Value = "'9.04"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.70%  '
$ws.Range('D29').Value = "'19.84"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.35%  '
$ws.Range('D30').Value = "'0.131"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -14.11%  '
$ws.Range('D31').Value = "'1.35"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.73%  '
$ws.Range('E33').Value = '  +5.97%  '
$ws.Range('D34').Value = "'4.67"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.21%  '
$ws.Range('D35').Value = "'4.50"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.77%  '
$ws.Range('D36').Value = "'2.46"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +6.35%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = "'3.42"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.00%  '
$ws.Range('B39').Value = 'WEMIXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D39').Value = "'1.82"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.17%  '
$ws.Range('D40').Value = "'5.30"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -4.39%  '
$ws.Range('E41').Value = '  +3.53%  '
$ws.Range('D42').Value = "'0.0963"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.43%  '
$ws.Range('D43').Value = "'0.0215"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.71%  '
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('D45').Value = "'1.401.29"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.14%  '
$ws.Range('D46').Value = "'16.02"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.89%  '
$ws.Range('D47').Value = "'90.10"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.31%  '
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('D49').Value = "'7.32"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('E50').Value = '  +1.64%  '
$ws.Range('D51').Value = "'2.02"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.44%  '
